# Generate Report for Handoff
#
# The "b.md" source file has completed a new handoff round in both the
# zh-cn and de-de localization sheets. Update:
#   - per-language Status / Latest Handoff File / Latest Handoff Datetime
#     for the b.md row (row 3) on the "zh-cn" and "de-de" sheets
#   - the hyperlink display text for the new handoff file name
#   - the roll-up Status / Latest Handoff Date on the "Overview" sheet

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

$zhFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhDatetime = "2016-03-19 07:54:50"

$deFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$deDatetime = "2016-03-19 07:54:59"

# --- zh-cn sheet: row 3 is the b.md row -------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = $newStatus
$zh.Range("D3").Value = $zhFile
$zh.Range("E3").Value = $zhDatetime

foreach ($h in $zh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$D$3') {
        $h.TextToDisplay = $zhFile
    }
}

# --- de-de sheet: row 3 is the b.md row -------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = $newStatus
$de.Range("D3").Value = $deFile
$de.Range("E3").Value = $deDatetime

foreach ($h in $de.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$D$3') {
        $h.TextToDisplay = $deFile
    }
}

# --- Overview sheet: roll up status + latest handoff date for b.md ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B3").Value = $newStatus
$ov.Range("C3").Value = $newStatus
$ov.Range("D3").Value = $deDatetime
